# Apply the DEU model update (2025-08-19 15:10) to the "wind" worksheet.
# The underlying change re-orders a handful of shared-string pairs for the
# won-DEU_40 / 34 / 33 / 25 / 21 / 19 wind-resource cost classes, which in turn
# causes the visible content (process name, description, cap_bnd, af~fx,
# ncap_cost, lcoe_class) of certain row pairs to swap places, while leaving
# cell formatting untouched. There is also a tiny floating point rounding
# refresh on row 41 (AF~FX / ncap_cost~USD21_alt).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("wind")

# Columns that carry the swapped data for each affected row pair:
#   C (process), D (description), K (comm-out "process" ref), M (cap_bnd),
#   N (af~fx), O (ncap_cost~USD21_alt), P (lcoe_class)
$cols = @("C", "D", "K", "M", "N", "O", "P")

function Swap-Rows($sheet, $rowA, $rowB, $columns) {
    foreach ($col in $columns) {
        $rangeA = $sheet.Range("$col$rowA")
        $rangeB = $sheet.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value = $valB
        $rangeB.Value = $valA
    }
}

$rowPairs = @(
    @(29, 30),
    @(53, 54),
    @(59, 60),
    @(99, 100),
    @(120, 121),
    @(127, 128)
)

foreach ($pair in $rowPairs) {
    Swap-Rows $ws $pair[0] $pair[1] $cols
}

# Minor floating point refresh on row 41 (values are numerically equal,
# just re-expressed with a slightly different binary representation).
$ws.Range("N41").Value = 0.3695
$ws.Range("O41").Value = 22.967898136436006
